$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalized counts (rows 1-23, columns B:G) replacing the raw Incucyte counts.
$rows = 1..23
$data = @(
    @(1608.368561735376, 1509.9051596532036, 1417.3751014610389, 1637.9784017278619, 1687.4419970631423, 1319.2566010058674),
    @(1642.7015523799737, 1530.1234933389721, 1436.8124999999998, 1661.9956803455725, 1731.2716593245225, 1389.5222129086335),
    @(1682.3165415852784, 1568.9426940156479, 1450.0299310064934, 1730.8452123830095, 1737.1718061674007, 1421.8917644593459),
    @(1765.9481854631442, 1633.6413618101078, 1511.4521103896102, 1786.8855291576676, 1773.4155653450805, 1498.4733864207878),
    @(1862.7848257427781, 1702.3836963417211, 1564.3218344155841, 1854.1339092872572, 1833.2599118942728, 1566.3704945515506),
    @(1952.5788012748026, 1801.0491647282722, 1663.8413149350647, 1950.2030237580996, 1930.1908957415562, 1660.3211441743501),
    @(2030.0481134985098, 1871.4089659547471, 1770.3582589285711, 1978.2231821454286, 1993.4067547723932, 1739.2712699077954),
    @(2102.2354271615095, 1928.8290336223301, 1820.1179991883114, 2087.9020878329734, 2065.894273127753, 1806.3788767812237),
    @(2199.0720674411436, 2001.6150348910974, 1923.5249594155841, 2149.5464362850976, 2148.4963289280468, 1860.8544635373007),
    @(2311.7547034028998, 2117.2639035736938, 2044.0368303571427, 2292.0489560835135, 2279.9853157121879, 1984.8061609388094),
    @(2531.8379767657043, 2297.6114400507504, 2200.3135146103891, 2414.5370770338377, 2483.1189427312775, 2124.5478834870073),
    @(2673.5716048113504, 2420.5389088602237, 2320.0478896103891, 2559.4413246940248, 2602.8076358296621, 2234.2885582564959),
    @(2864.603886090265, 2593.6078452104039, 2463.884638798701, 2721.1576673866093, 2762.1116005873714, 2414.294844928751),
    @(3045.952503341216, 2764.2505815182913, 2656.7036323051943, 2888.4780417566599, 2886.8575624082227, 2580.0901089689855),
    @(3189.4467975737643, 2912.2487840981175, 2854.9650974025972, 3043.789776817855, 3060.4904552129219, 2721.4108340318521),
    @(3409.5300709365692, 3039.219919644745, 2956.0395698051943, 3175.884809215263, 3181.0220264317177, 2884.0480930427489),
    @(3575.9130255988493, 3155.6775216747724, 3125.5336850649346, 3354.4132469402452, 3288.0675477239351, 3027.7373218776193),
    @(3715.8859874575933, 3285.8835906111226, 3248.3780438311683, 3437.6731461483082, 3376.5697503671067, 3153.2680217937968),
    @(3810.0816284568737, 3396.6800592091349, 3344.0100446428569, 3532.1411087113033, 3512.2731277533035, 3241.6921626152553),
    @(3884.9099414002271, 3454.9088602241486, 3455.1919642857138, 3572.9704823614115, 3566.217327459618, 3363.2753562447606),
    @(3994.9515780816291, 3531.7385282300693, 3515.0591517857138, 3651.4269258459326, 3648.8193832599113, 3434.3304694048611),
    @(4124.3605428189585, 3623.1253964897437, 3644.9009740259735, 3770.7127429805619, 3825.8237885462549, 3547.2291492036879),
    @(4229.1201809396534, 3757.3751321632476, 3769.3003246753242, 3851.5709143268541, 3953.9412628487512, 3694.0763830678957)
)

for ($idx = 0; $idx -lt $rows.Length; $idx++) {
    $r = $rows[$idx]
    $vals = $data[$idx]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

# Mirror the author's final selection (B1:G23) as recorded in the saved workbook.
$ws.Range("B1:G23").Select() | Out-Null
